$d = $word.ActiveDocument

# Locate the title text " de coador." so we don't depend on hard-coded offsets.
$found = $d.Content
$ok = $found.Find.Execute(" de coador.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($ok) {
    # The trailing period is the very last character of the match.
    $periodRange = $d.Range($found.End - 1, $found.End)

    # Insert the "_GoBack" bookmark collapsed right before that period
    # (i.e. right after "coador"). Bookmarks.Add re-uses the existing
    # bookmark name, so this effectively moves the bookmark that used to
    # sit in the empty trailing paragraph up to this spot.
    $bookmarkRange = $d.Range($found.End - 1, $found.End - 1)
    $d.Bookmarks.Add("_GoBack", $bookmarkRange)

    # Now drop the trailing period itself.
    $periodRange.Delete()
}
